$wb = $excel.ActiveWorkbook

# --- Status update: "Ready for handoff" -> "In Translation" -------------
# Overview sheet tracks status per-language in columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# Per-locale sheets track status in column C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the Status columns -------------------------------------------
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
